# Updates the cryptos list: refresh Price (D) and Volume(1h) (E) columns,
# and correct the ranking swap between OKB and FirstDigitalUSD (rows 42/43).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 42/43 swapped (OKB now ranks above FirstDigitalUSD) ---
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "'36.41"
$ws.Range("E42").Value = "  -9.50%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  -0.03%  "

# --- Refreshed Price / Volume(1h) values for the remaining rows ---
$ws.Range("D2").Value = "55.655.80"
$ws.Range("E2").Value = "  -2.92%  "
$ws.Range("D3").Value = "2.955.05"
$ws.Range("E3").Value = "  -5.19%  "
$ws.Range("D5").Value = "'491.10"
$ws.Range("E5").Value = "  -5.51%  "
$ws.Range("D6").Value = "'132.84"
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "2.956.38"
$ws.Range("E8").Value = "  -5.03%  "
$ws.Range("D9").Value = "'0.422"
$ws.Range("E9").Value = "  -5.08%  "
$ws.Range("D10").Value = "'7.13"
$ws.Range("E10").Value = "  -0.39%  "
$ws.Range("E11").Value = "  -6.09%  "
$ws.Range("D12").Value = "'0.347"
$ws.Range("E12").Value = "  -8.46%  "
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("D14").Value = "3.459.46"
$ws.Range("E14").Value = "  -5.21%  "
$ws.Range("D15").Value = "'24.67"
$ws.Range("E15").Value = "  -2.32%  "
$ws.Range("D16").Value = "55.583.20"
$ws.Range("E16").Value = "  -3.02%  "
$ws.Range("D17").Value = "2.946.92"
$ws.Range("E17").Value = "  -5.40%  "
$ws.Range("D18").Value = "'0.0000141"
$ws.Range("E18").Value = "  -5.13%  "
$ws.Range("D19").Value = "'5.65"
$ws.Range("D20").Value = "'12.15"
$ws.Range("E20").Value = "  -5.56%  "
$ws.Range("D21").Value = "'7.58"
$ws.Range("E21").Value = "  -4.71%  "
$ws.Range("D22").Value = "'317.88"
$ws.Range("E22").Value = "  -7.09%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'0.461"
$ws.Range("E24").Value = "  -8.15%  "
$ws.Range("D25").Value = "'60.13"
$ws.Range("E25").Value = "  -11.90%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.61%  "
$ws.Range("D27").Value = "'0.160"
$ws.Range("E27").Value = "  -3.11%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").Value = "0.0₃0846"
$ws.Range("E29").Value = "  -9.06%  "
$ws.Range("D30").Value = "'6.46"
$ws.Range("E30").Value = "  -3.88%  "
$ws.Range("D31").Value = "'6.56"
$ws.Range("E31").Value = "  -5.39%  "
$ws.Range("D33").Value = "'1.70"
$ws.Range("E33").Value = "  -8.25%  "
$ws.Range("D34").Value = "'19.49"
$ws.Range("E34").Value = "  -8.86%  "
$ws.Range("D35").Value = "'149.17"
$ws.Range("E35").Value = "  -6.09%  "
$ws.Range("D36").Value = "'4.38"
$ws.Range("E36").Value = "  -7.76%  "
$ws.Range("D37").Value = "'1.31"
$ws.Range("E37").Value = "  -4.78%  "
$ws.Range("D38").Value = "'5.68"
$ws.Range("E38").Value = "  -7.29%  "
$ws.Range("D39").Value = "'0.0656"
$ws.Range("E39").Value = "  -4.00%  "
$ws.Range("D40").Value = "'23.34"
$ws.Range("E40").Value = "  -6.81%  "
$ws.Range("D41").Value = "2.985.95"
$ws.Range("E41").Value = "  -4.99%  "
$ws.Range("D48").Value = "2.114.48"
$ws.Range("E48").Value = "  -6.33%  "
$ws.Range("E49").Value = "  +1.15%  "
$ws.Range("D50").Value = "'19.28"
$ws.Range("E50").Value = "  -2.00%  "
$ws.Range("D51").Value = "'5.54"
$ws.Range("E51").Value = "  -9.79%  "

Write-Output "cryptos list updated"
